$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.078.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.22%  "

$ws.Range("D3").Value = "'3.546.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +9.27%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "'190.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.09%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'566.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.55%  "

$ws.Range("D7").Value = "'3.540.65"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.15%  "

$ws.Range("D8").Value = "'0.619"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.83%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").Value = "'0.636"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.61%  "

$ws.Range("D11").Value = "'0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.55%  "

$ws.Range("D12").Value = "'55.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.48%  "

$ws.Range("D13").Value = "'0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.42%  "

$ws.Range("D14").Value = "'9.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.28%  "

$ws.Range("D15").Value = "'4.101.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.24%  "

$ws.Range("D16").Value = "'3.542.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.39%  "

$ws.Range("E17").Value = "  +4.04%  "

$ws.Range("D18").Value = "'67.047.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.37%  "

$ws.Range("D19").Value = "'18.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.18%  "

$ws.Range("D20").Value = "'12.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.62%  "

$ws.Range("D21").Value = "'1.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.63%  "

$ws.Range("D22").Value = "'434.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +17.80%  "

$ws.Range("E23").Value = "  +10.61%  "

$ws.Range("D24").Value = "'85.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.53%  "

$ws.Range("D25").Value = "'4.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.78%  "

$ws.Range("D26").Value = "'11.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.68%  "

$ws.Range("D27").Value = "'2.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.61%  "

$ws.Range("D28").Value = "'12.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.75%  "

$ws.Range("D29").Value = "'9.18"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").Value = "'30.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.21%  "

$ws.Range("D31").Value = "'642.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("D32").Value = "'6.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.63%  "

$ws.Range("D33").Value = "'11.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.84%  "

$ws.Range("E34").Value = "  +5.21%  "

$ws.Range("D35").Value = "'59.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.45%  "

$ws.Range("D36").Value = "'38.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.58%  "

$ws.Range("D37").Value = "'0.0₃0816"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +12.26%  "

$ws.Range("E38").Value = "  +18.73%  "

$ws.Range("E39").Value = "  -0.15%  "

$ws.Range("D40").Value = "'0.392"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.93%  "

$ws.Range("D41").Value = "'3.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.32%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("D43").Value = "'3.043.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.70%  "

$ws.Range("D44").Value = "'2.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.04%  "

$ws.Range("E45").Value = "  +11.19%  "

$ws.Range("D46").Value = "'3.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.63%  "

$ws.Range("D47").Value = "'0.0420"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.19%  "

$ws.Range("D48").Value = "'2.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.54%  "

$ws.Range("E49").Value = "  +5.94%  "

$ws.Range("D50").Value = "'142.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.80%  "

$ws.Range("D51").Value = "'8.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.21%  "
